$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    ,@(2, -0.2103977471099733, -0.3051268707526743)
    ,@(3, -0.4647572337589958, -0.3974561397256904)
    ,@(4, -0.5420458476225329, -0.4774627561395389)
    ,@(5, -0.2020935448562471, -0.241371250993697)
    ,@(6, -0.1988301985982694, -0.1807974864889489)
    ,@(7, -0.2383226232223296, -0.2892530252971715)
    ,@(8, -0.479677295720217, -0.414986518278334)
    ,@(9, -0.6423532367102, -0.5017630512574407)
    ,@(10, -0.4975737005612251, -0.2802299313223474)
    ,@(11, -0.2691207381554653, -0.2961359379524409)
    ,@(12, -0.2204045929313549, -0.2012800297013035)
    ,@(13, -0.1018331133741392, -0.08547658697900877)
    ,@(14, -0.2709067696138323, -0.2243252943307559)
    ,@(15, -0.06910856589356343, -0.02463965958345737)
    ,@(16, -0.2385271655087465, -0.0995548935257257)
    ,@(17, 0.003973537876794532, 0.1244900861441573)
    ,@(18, 0.02150724881991195, 0.05279657228616021)
    ,@(19, 0.05510132201968172, 0.1004400509282022)
    ,@(20, -0.1575318625013999, -0.09097624358877618)
    ,@(21, -0.03319693324335826, 0.03147393684190031)
    ,@(22, -0.07163056772763338, 0.07302722220851199)
    ,@(23, 0.02687111395900857, 0.02539431603911637)
    ,@(24, 0.8395935783534152, 0.6764159873527081)
    ,@(25, 0.1298358672208684, 0.1100388389983356)
    ,@(26, 0.1838494099857298, 0.1536557422819028)
    ,@(27, 0.1310435795481181, 0.1346133361584621)
    ,@(28, 0.3703360383785936, 0.2314731374672672)
    ,@(29, 0.6984110970134783, 0.5849379420876908)
    ,@(30, 0.2505934258380249, 0.2025572122821401)
    ,@(31, 0.08202811042531011, 0.0652986017292802)
    ,@(32, 0.2320925800028873, 0.2164869098289154)
    ,@(33, 0.1518265953666065, 0.1599713753107591)
    ,@(34, 0.1094245586281686, 0.07668897726859338)
    ,@(35, 0.4042957538170733, 0.2904597654567493)
    ,@(36, 0.2842324891237176, 0.1613062843236384)
    ,@(37, -0.02671313311160461, -0.01921542017419373)
    ,@(38, 0.3607894319027811, 0.3034417197658699)
    ,@(39, -0.1124807735437134, -0.1498292724577766)
    ,@(40, 0.1492292611282714, 0.183457992883758)
    ,@(41, -0.1645398281152756, -0.2194875342870689)
    ,@(42, 0.2819136431339014, 0.2764821163800744)
    ,@(43, 0.1166719301517765, 0.1570997822562603)
    ,@(44, 0.03903648804324042, -0.03261894322829163)
    ,@(45, -0.1207317465490824, -0.07712475407900811)
    ,@(46, -0.172434941020288, -0.1607961516556028)
    ,@(47, -0.1679883498925968, -0.1598627741767813)
    ,@(48, -0.2209854448837378, -0.1935945540093548)
    ,@(49, -0.2234889888131208, -0.1933404368588273)
    ,@(50, -0.1022061371806729, -0.1135780396913113)
    ,@(51, -0.2278083440056071, -0.2170910917395659)
    ,@(52, -0.2278083440056071, -0.2170910917395659)
    ,@(53, -0.1952448759629271, -0.1756372678191264)
    ,@(54, -0.2180694839001492, -0.1926282370764213)
    ,@(55, -0.1552409275042511, -0.1361822745476835)
    ,@(56, -0.1175524768941751, -0.1160085596124987)
    ,@(57, -0.1837436120102373, -0.1421331029364994)
    ,@(58, -0.1467449462786686, -0.1967644572964131)
    ,@(59, -0.2224792430287997, -0.2299764227862012)
    ,@(60, -0.2531512652069413, -0.24948062307687)
    ,@(61, -0.2679013037651269, -0.2248580835488676)
    ,@(62, -0.1498573848783344, -0.09214075319993889)
    ,@(63, -0.4712532092516, -0.4407478773075003)
    ,@(64, -0.2960589708242725, -0.2749834954519098)
    ,@(65, -0.2482426222997058, -0.2384391553372095)
    ,@(66, -0.1047743062095829, -0.08620917232219795)
    ,@(67, 0.01947647773925388, -0.01456983493393927)
)

foreach ($entry in $newValues) {
    $r = $entry[0]
    $ws.Cells.Item($r, 1).Value = $entry[1]
    $ws.Cells.Item($r, 2).Value = $entry[2]
}

Write-Host "Updated rows 2-67 with new standardized PCA values"
